# Commit: "Fri, Jul 31, 2020  9:06:23 AM"
#
# The deck's Design/Theme is switched from the custom "Integral" theme to
# the stock "Office Theme" (i.e. the Office default design was (re-)applied
# from the Design gallery). In the canonical OOXML this shows up as:
#   - ppt/theme/theme1.xml (the theme bound to the slide master, i.e. the
#     deck's working theme) gets the "Office Theme" color palette, while
#   - the table on slide 6 gets re-stamped with the built-in "medium" table
#     style GUID that Office assigns to tables when the active theme no
#     longer matches the table's previous (theme-specific) style id.
#
# ThemeColorScheme writes always land on the deck's live theme part, so we
# drive the 12 theme colors to the stock Office palette, and re-apply the
# table style via Table.ApplyStyle (the only writable entry point for a
# table's style id — Table.Style.Id/Table.Style are read-only / reject
# direct assignment).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Swap the working theme's color scheme from "Integral" to the stock
#    "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in the
#    order exposed by ThemeColorScheme — OLE RGB longs, i.e. 0xBBGGRR).
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0         # dk1      000000
$scheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388   # dk2      44546A
$scheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407     # accent4  FFC000
$scheme.Item(9).RGB  = 12874308  # accent5  4472C4
$scheme.Item(10).RGB = 4697456   # accent6  70AD47
$scheme.Item(11).RGB = 12673797  # hlink    0563C1
$scheme.Item(12).RGB = 7491477   # folHlink 954F72

# ---------------------------------------------------------------------
# 2) Re-stamp the table on slide 6 with the built-in table style that
#    Office uses once the theme-specific style no longer applies.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{5BCFD7C7-0DDE-4093-B28B-2102F57E609D}")
    }
}
